$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a 20-row data block (rows 2-21) that already repeats
# six times (rows 2-121). Append two more copies of that block so the
# data spans rows 2-161 (dimension A1:N161).
$block = $ws.Range("A2:N21")
$block.Copy($ws.Range("A122"))
$block.Copy($ws.Range("A142"))
